# Bento "startup" tab: refresh the grouped_recurrence_score filter used by
# the Neo4j/db queries from "51-100" down to "0-5" (new data set availability).
# Row 2 = CasesTab, Row 3 = SamplesTab, Row 4 = FilesTab.
# Column C holds the shared StatQuery (count) query on every row; column B
# holds each tab's own db query. B2 (CasesTab's own query happens to live in
# a different slot and is intentionally left untouched, matching the source
# edit) keeps its original "51-100" filter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldFilter = '51-100'
$newFilter = '0-5'
$pattern = [regex]::Escape($oldFilter)

$cellsToUpdate = @('C2', 'B3', 'C3', 'B4', 'C4')
foreach ($addr in $cellsToUpdate) {
    $rng = $ws.Range($addr)
    $current = $rng.Value2
    if ($current -ne $null -and $current -match $pattern) {
        $rng.Value2 = ($current -replace $pattern, $newFilter)
    }
}

# Restore the view focus to the top data row / its StatQuery cell, matching
# the saved workbook's scroll + selection state.
$ws.Application.ActiveWindow.ScrollRow = 2
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range('C2').Select()
